# Update countries & provincias Spain
#
# Refreshes the COVID-19 "Pais" sheet with the next data snapshot:
#  - Updates the "Datos actualizados a ..." timestamp cell (A1).
#  - Writes the latest totals (Casos totales, Nuevos casos, Casos activos,
#    Recuperados, Casos criticos, Muertes hoy, Muertes) for every country
#    row whose figures moved.
#  - A few countries leap-frogged their neighbours in the ranking
#    (Sudafrica/Arabia Saudita, Etiopia/Haiti/Tayikistan/Bulgaria/
#    Gabon/Bosnia y Herzegovina/Costa Rica, Cabo Verde/Nueva Zelanda,
#    Zimbabue/Republica del Chad/Principado de Andorra,
#    Groenlandia/Islas Malvinas); those rows get both a new country name
#    (column A) and new figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 -> Datos actualizados a 9 de Julio de 2020 a las 00:18
$ws.Range("A1").Value = 'Datos actualizados a 9 de Julio de 2020 a las 00:18'

# Row 4 -> Estados Unidos
$ws.Range("B4").Value = 3148275
$ws.Range("C4").Value = 51191
$ws.Range("D4").Value = 1384622
$ws.Range("E4").Value = 1628918
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 763
$ws.Range("H4").Value = 134735

# Row 5 -> Brasil
$ws.Range("B5").Value = 1713160
$ws.Range("C5").Value = 38505
$ws.Range("D5").Value = 1117922
$ws.Range("E5").Value = 527274
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1096
$ws.Range("H5").Value = 67964

# Row 16 -> Sudafrica
$ws.Range("A16").Value = 'Sudafrica'
$ws.Range("B16").Value = 224665
$ws.Range("C16").Value = 8810
$ws.Range("D16").Value = 106842
$ws.Range("E16").Value = 114223
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 98
$ws.Range("H16").Value = 3600

# Row 17 -> Arabia Saudita
$ws.Range("A17").Value = 'Arabia Saudita'
$ws.Range("B17").Value = 220144
$ws.Range("C17").Value = 3036
$ws.Range("D17").Value = 158050
$ws.Range("E17").Value = 60035
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 42
$ws.Range("H17").Value = 2059

# Row 19 -> Alemania
$ws.Range("B19").Value = 198765
$ws.Range("C19").Value = 410
$ws.Range("D19").Value = 182700
$ws.Range("E19").Value = 6950
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 12
$ws.Range("H19").Value = 9115

# Row 50 -> Barein
$ws.Range("B50").Value = 30931
$ws.Range("C50").Value = 610
$ws.Range("D50").Value = 26073
$ws.Range("E50").Value = 4760
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 98

# Row 84 -> Etiopia
$ws.Range("A84").Value = 'Etiopia'
$ws.Range("B84").Value = 6774
$ws.Range("C84").Value = 928
$ws.Range("D84").Value = 2430
$ws.Range("E84").Value = 4224
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 17
$ws.Range("H84").Value = 120

# Row 85 -> Haiti
$ws.Range("A85").Value = 'Haiti'
$ws.Range("B85").Value = 6432
$ws.Range("C85").Value = 61
$ws.Range("D85").Value = 2080
$ws.Range("E85").Value = 4235
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 4
$ws.Range("H85").Value = 117

# Row 86 -> Tayikistan
$ws.Range("A86").Value = 'Tayikistan'
$ws.Range("B86").Value = 6364
$ws.Range("C86").Value = 49
$ws.Range("D86").Value = 5011
$ws.Range("E86").Value = 1299
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 54

# Row 87 -> Bulgaria
$ws.Range("A87").Value = 'Bulgaria'
$ws.Range("B87").Value = 6342
$ws.Range("C87").Value = 240
$ws.Range("D87").Value = 3166
$ws.Range("E87").Value = 2917
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 5
$ws.Range("H87").Value = 259

# Row 88 -> Gabon
$ws.Range("A88").Value = 'Gabon'
$ws.Range("B88").Value = 5871
$ws.Range("C88").Value = 128
$ws.Range("D88").Value = 2682
$ws.Range("E88").Value = 3143
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 46

# Row 89 -> Bosnia y Herzegovina
$ws.Range("A89").Value = 'Bosnia y Herzegovina'
$ws.Range("B89").Value = 5869
$ws.Range("C89").Value = 248
$ws.Range("D89").Value = 2769
$ws.Range("E89").Value = 2891
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 2
$ws.Range("H89").Value = 209

# Row 90 -> Costa Rica
$ws.Range("A90").Value = 'Costa Rica'
$ws.Range("B90").Value = 5836
$ws.Range("C90").Value = 350
$ws.Range("D90").Value = 1929
$ws.Range("E90").Value = 3883
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 24

# Row 125 -> Cabo Verde
$ws.Range("A125").Value = 'Cabo Verde'
$ws.Range("B125").Value = 1542
$ws.Range("C125").Value = 43
$ws.Range("D125").Value = 730
$ws.Range("E125").Value = 794
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 18

# Row 126 -> Nueva Zelanda
$ws.Range("A126").Value = 'Nueva Zelanda'
$ws.Range("B126").Value = 1537
$ws.Range("C126").Value = 1
$ws.Range("D126").Value = 1492
$ws.Range("E126").Value = 23
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 22

# Row 145 -> Zimbabue
$ws.Range("A145").Value = 'Zimbabue'
$ws.Range("B145").Value = 885
$ws.Range("C145").Value = 98
$ws.Range("D145").Value = 206
$ws.Range("E145").Value = 670
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 9

# Row 146 -> Republica del Chad
$ws.Range("A146").Value = 'Republica del Chad'
$ws.Range("B146").Value = 873
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 788
$ws.Range("E146").Value = 11
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 74

# Row 147 -> Principado de Andorra
$ws.Range("A147").Value = 'Principado de Andorra'
$ws.Range("B147").Value = 855
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 800
$ws.Range("E147").Value = 3
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 52

# Row 149 -> Santo Tome y Principe
$ws.Range("B149").Value = 724
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 283
$ws.Range("E149").Value = 428
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 13

# Row 209 -> Groenlandia
$ws.Range("A209").Value = 'Groenlandia'

# Row 210 -> Islas Malvinas
$ws.Range("A210").Value = 'Islas Malvinas'
